# Generate Report for Handoff
# Adds a new "d67eb8db-77e3-4913-a2b6-f8b83e1b5f18" handoff row (row 3) to each
# of the three worksheets (Overview, zh-cn, de-de), mirroring the existing
# "2cc5c88e-c09b-468f-b7a1-ba9abf434523" row already present in row 2.

$wb = $excel.ActiveWorkbook

$mdName   = "d67eb8db-77e3-4913-a2b6-f8b83e1b5f18.md"
$zhXlf    = "d67eb8db-77e3-4913-a2b6-f8b83e1b5f18.ae00ddcaf782e69bff1ee3e5f4a3b8c8062299ee.zh-cn.xlf"
$deXlf    = "d67eb8db-77e3-4913-a2b6-f8b83e1b5f18.ae00ddcaf782e69bff1ee3e5f4a3b8c8062299ee.de-de.xlf"

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/dd042ec8035629dc3e0850a7ce49ff25f31bcae9/e2e/$mdName"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/867726822f925ecc5851d070c1a642b1352510eb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e9e313155e43d205acb198100a17c5f960e93b3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value2 = $mdName
$wsOverview.Range("B3").Value2 = "Ready for handoff"
$wsOverview.Range("C3").Value2 = "Ready for handoff"
$wsOverview.Range("D3").Value2 = "2016-30-20 20:30:28"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, "", "", $mdName)

# ---------------------------------------------------------------------------
# zh-cn sheet: Source File Name | File Extension | Status | Latest Handoff File |
#              Latest Handoff Datetime | ... | Latest Handback DateTime | Handoff Reason
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value2 = $mdName
$wsZh.Range("B3").Value2 = ".md"
$wsZh.Range("C3").Value2 = "Ready for handoff"
$wsZh.Range("D3").Value2 = $zhXlf
$wsZh.Range("E3").Value2 = "2016-03-20 20:30:25"
$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").Value2 = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value2 = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, "", "", $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $mdUrl, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhXlfUrl, "", "", $zhXlf)

# ---------------------------------------------------------------------------
# de-de sheet: same shape as zh-cn
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value2 = $mdName
$wsDe.Range("B3").Value2 = ".md"
$wsDe.Range("C3").Value2 = "Ready for handoff"
$wsDe.Range("D3").Value2 = $deXlf
$wsDe.Range("E3").Value2 = "2016-03-20 20:30:28"
$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").Value2 = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value2 = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, "", "", $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $mdUrl, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deXlfUrl, "", "", $deXlf)

Write-Output "Handoff row added for $mdName"
